$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "week" data: complete the 20/05 pair, add 27/05, stub 3/06 and the
#     final 7/06 TO-DO, then fill in the 27/05 TO-DO details. Filled in the
#     order a real editing session would naturally follow so new shared
#     strings are appended in the same sequence as the target workbook.
$ws.Range("G2").Value = "Done 20/05/2021"
$ws.Range("H2").Value = "TO-DO 27/05/2021"

$ws.Range("G3").Value = "Art-work Q*bert"

$ws.Range("I2").Value = "Done 27/05/2021"
$ws.Range("J2").Value = "TO-DO 3/06/2021"
$ws.Range("K2").Value = "Done 3/06/2021"
$ws.Range("L2").Value = "TO-DO Final 7/06/2021"

$ws.Range("H3").Value = "Game-board Implementeren 4 slots"
$ws.Range("H4").Value = "Q*bert player char 2 slots"
$ws.Range("H5").Value = "Teleportation disks 2 slots"
$ws.Range("H6").Value = "Coily enemy char 4 slots"
$ws.Range("H7").Value = "Point system 2 slots"

# --- Column widths for the newly introduced columns G:L ---
$ws.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws.Columns.Item(8).ColumnWidth = 32.333333333333336
$ws.Columns.Item(9).ColumnWidth = 15.166666666666666
$ws.Columns.Item(10).ColumnWidth = 15.333333333333334
$ws.Columns.Item(11).ColumnWidth = 14.166666666666666
$ws.Columns.Item(12).ColumnWidth = 20.166666666666668

# --- Selection moves to H9 (matches the saved view state) ---
$ws.Range("H9").Select() | Out-Null
